# -----------------------------------------------------------------------
# Livestock module update:
#  - insert new "adm_src_*" parameter rows into the "conf" sheet
#  - rename the "plot_character_file" parameter to "plo_character_file"
#  - append new "plo_src_*" / "plo_lim_*" parameter rows
#  - add a new "buffer" sheet that will hold the per-plot region layer
# -----------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$conf = $wb.Worksheets.Item("conf")

# ---------------------------------------------------------------------
# 1) Insert 8 new blank rows above the old row 10 ("plot_character_file")
#    This pushes that row down to row 18 and gives us rows 10-17 to fill
#    in with the new adm_src_* parameters.
# ---------------------------------------------------------------------
$conf.Rows("10:17").Insert()

# Newly inserted rows inherit the hyperlink formatting of the row above
# (row 9 has a hyperlink). Clear that back out before restyling.
$conf.Range("A10:B17").ClearFormats()

function Set-ConfRow($row, $paramName, $paramValue) {
    $a = $conf.Cells.Item($row, 1)
    $b = $conf.Cells.Item($row, 2)
    $a.Value = $paramName
    $b.Value = $paramValue
    $b.Font.Name = "Calibri"
    $b.Font.Size = 11
}

Set-ConfRow 10 "adm_src_adm1_name" "NOM_DEP"
Set-ConfRow 11 "adm_src_adm1_id"   "COD_DPTO"
Set-ConfRow 12 "adm_src_adm2_name" "NOMB_MPIO"
Set-ConfRow 13 "adm_src_adm2_id"   "DPTOMPIO"
Set-ConfRow 14 "adm_src_adm3_name" "NOMBRE_VER"
Set-ConfRow 15 "adm_src_adm3_id"   "CODIGO_VER"
Set-ConfRow 16 "adm_src_area"      "AREA_HA"
Set-ConfRow 17 "adm_src_geometry"  "geometry"

# ---------------------------------------------------------------------
# 2) Row 18 is the old "plot_character_file" row (pushed down by the
#    insert above). Fix the typo in the parameter name; the value (";")
#    is unchanged.
# ---------------------------------------------------------------------
$conf.Cells.Item(18, 1).Value = "plo_character_file"

# ---------------------------------------------------------------------
# 3) Append the new plo_src_* / plo_lim_* parameter rows after row 18.
# ---------------------------------------------------------------------
$conf.Cells.Item(19, 1).Value = "plo_src_ext_id"
$conf.Cells.Item(19, 2).Value = "codigosit"

$conf.Cells.Item(20, 1).Value = "plo_src_lat"
$conf.Cells.Item(20, 2).Value = "latitud"

$conf.Cells.Item(21, 1).Value = "plo_src_lon"
$conf.Cells.Item(21, 2).Value = "longitud"

$conf.Cells.Item(22, 1).Value = "plo_src_animals"
$conf.Cells.Item(22, 2).Value = "totalanimales"

$conf.Cells.Item(23, 1).Value = "plo_src_adm1"
$conf.Cells.Item(23, 2).Value = "departamento"

$conf.Cells.Item(24, 1).Value = "plo_src_adm2"
$conf.Cells.Item(24, 2).Value = "municipio"

$conf.Cells.Item(25, 1).Value = "plo_src_adm3"
$conf.Cells.Item(25, 2).Value = "vereda"

$conf.Cells.Item(26, 1).Value = "plo_src_crs"
$conf.Cells.Item(26, 2).Value = 4326

$conf.Cells.Item(27, 1).Value = "plo_lim_lat_low"
$conf.Cells.Item(27, 2).Value = -4.3

$conf.Cells.Item(28, 1).Value = "plo_lim_lat_upp"
$conf.Cells.Item(28, 2).Value = 12.4

$conf.Cells.Item(29, 1).Value = "plo_lim_lon_low"
$conf.Cells.Item(29, 2).Value = -79

$conf.Cells.Item(30, 1).Value = "plo_lim_lon_upp"
$conf.Cells.Item(30, 2).Value = -66

# ---------------------------------------------------------------------
# 4) Update the "conf" sheet view: scroll down and select A27:A30.
# ---------------------------------------------------------------------
$conf.Range("A27:A30").Select()
$excel.ActiveWindow.ScrollRow = 15

# ---------------------------------------------------------------------
# 5) Add a new "buffer" sheet after "conf" for the per-plot region /
#    field-capacity lookup layer.
# ---------------------------------------------------------------------
$buffer = $wb.Worksheets.Add([System.Type]::Missing, $conf)
$buffer.Name = "buffer"

$buffer.Cells.Item(1, 1).Value = "region_id"
$buffer.Cells.Item(1, 2).Value = "region"
$buffer.Cells.Item(1, 3).Value = "field_capacity"

$buffer.Columns.Item(1).ColumnWidth = 8.5546875
$buffer.Columns.Item(2).ColumnWidth = 11
$buffer.Columns.Item(3).ColumnWidth = 12.33203125

$buffer.Range("A2").Select()

$wb.Save()
